# Add data for 2024-11-05.
#
# The sheet has one date-labeled column per day; the new day's data is
# appended as a new rightmost column. In the source edit the existing last
# column ("BF", 2024/11/06) formatting/values are duplicated into the new
# column "BG", the header in BG is set back to the old "2024/11/06" date,
# and BF's header becomes the new "2024/11/05" date - exactly mirroring the
# target OOXML diff (numeric rows 2-53 simply carry BF's value over into
# the new BG cell; only the header row's text actually differs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Clone column BF's formatting (fill colors / styles) into the new
#    column BG for every row in one shot, so each row's BG cell ends up
#    with the same style (s=1/2/3) as its BF cell.
$ws.Range("BF1:BF53").Copy()
$ws.Range("BG1:BG53").PasteSpecial(-4122)

# 2) Give column BG the same column width as the rest of the data columns.
$ws.Columns("BG").ColumnWidth = 11.17

# 3) Header row: BG1 takes over the old "2024/11/06" label that used to
#    live in BF1. Force text (not auto-parsed date) via a temporary "@"
#    number format, then restore BF1's original style over it.
$ws.Range("BG1").NumberFormat = "@"
$ws.Cells.Item(1, 59).Value = "2024/11/06"
$ws.Range("BF1").Copy()
$ws.Range("BG1").PasteSpecial(-4122)

# BF1 becomes the new "2024/11/05" label - same trick, restoring style
# from the untouched BE1 header cell afterwards.
$ws.Range("BF1").NumberFormat = "@"
$ws.Cells.Item(1, 58).Value = "2024/11/05"
$ws.Range("BE1").Copy()
$ws.Range("BF1").PasteSpecial(-4122)

# 4) Rows 2-53: duplicate BF's numeric value into the new BG cell.
$ws.Cells.Item(2, 59).Value = 189.6
$ws.Cells.Item(3, 59).Value = 205.9
$ws.Cells.Item(4, 59).Value = 164.3
$ws.Cells.Item(5, 59).Value = 136.3
$ws.Cells.Item(6, 59).Value = 140.4
$ws.Cells.Item(7, 59).Value = 135.2
$ws.Cells.Item(8, 59).Value = 149.4
$ws.Cells.Item(9, 59).Value = 157.3
$ws.Cells.Item(10, 59).Value = 121
$ws.Cells.Item(11, 59).Value = 144
$ws.Cells.Item(12, 59).Value = 156.5
$ws.Cells.Item(13, 59).Value = 160.7
$ws.Cells.Item(14, 59).Value = 119.6
$ws.Cells.Item(15, 59).Value = 133.7
$ws.Cells.Item(16, 59).Value = 149.9
$ws.Cells.Item(17, 59).Value = 165.7
$ws.Cells.Item(18, 59).Value = 142.5
$ws.Cells.Item(19, 59).Value = 115.7
$ws.Cells.Item(20, 59).Value = 180.1
$ws.Cells.Item(21, 59).Value = 158.1
$ws.Cells.Item(22, 59).Value = 121.2
$ws.Cells.Item(23, 59).Value = 168.4
$ws.Cells.Item(24, 59).Value = 145.3
$ws.Cells.Item(25, 59).Value = 140.8
$ws.Cells.Item(26, 59).Value = 124.5
$ws.Cells.Item(27, 59).Value = 153.4
$ws.Cells.Item(28, 59).Value = 117
$ws.Cells.Item(29, 59).Value = 205.1
$ws.Cells.Item(30, 59).Value = 124.7
$ws.Cells.Item(31, 59).Value = 156.5
$ws.Cells.Item(32, 59).Value = 148.5
$ws.Cells.Item(33, 59).Value = 130.5
$ws.Cells.Item(34, 59).Value = 165.4
$ws.Cells.Item(35, 59).Value = 146.7
$ws.Cells.Item(36, 59).Value = 135.7
$ws.Cells.Item(37, 59).Value = 167.4
$ws.Cells.Item(38, 59).Value = 127.3
$ws.Cells.Item(39, 59).Value = 133.3
$ws.Cells.Item(40, 59).Value = 159.5
$ws.Cells.Item(41, 59).Value = 136
$ws.Cells.Item(42, 59).Value = 199
$ws.Cells.Item(43, 59).Value = 118.9
$ws.Cells.Item(44, 59).Value = 168.8
$ws.Cells.Item(45, 59).Value = 103.1
$ws.Cells.Item(46, 59).Value = 127.6
$ws.Cells.Item(47, 59).Value = 151.1
$ws.Cells.Item(48, 59).Value = 154.4
$ws.Cells.Item(49, 59).Value = 127.6
$ws.Cells.Item(50, 59).Value = 150.9
$ws.Cells.Item(51, 59).Value = 108.6
$ws.Cells.Item(52, 59).Value = 138.8
$ws.Cells.Item(53, 59).Value = 155.9
